$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row ("Row Index" -> "Index", "Product Name" -> "Product Detail")
$ws.Range("A1").Value = "Index"
$ws.Range("B1").Value = "Product Detail"

# Update product name (col B) and price (col C) for each data row (rows 2-31)
$ws.Range("B2").Value = "Meridian Round Solid Wood Coffee Table In Teak Finish"
$ws.Range("C2").Value = "₹5,939"
$ws.Range("B3").Value = "Adele Rectangular Engineered Wood Coffee Table In Classic Walnut Finish"
$ws.Range("C3").Value = "₹2,927"
$ws.Range("B4").Value = "Awdry Rectangular Engineered Wood Coffee Table In Sonoma Oak Finish"
$ws.Range("C4").Value = "₹2,903"
$ws.Range("B5").Value = "Altura Rectangular Solid Wood Coffee Table In Two Tone Finish"
$ws.Range("C5").Value = "₹14,453"
$ws.Range("B6").Value = "Claire Rectangular Solid Wood Coffee Table In Teak Finish"
$ws.Range("C6").Value = "₹12,725"
$ws.Range("B7").Value = "Tate Square Solid Wood Coffee Table In Teak Finish"
$ws.Range("C7").Value = "₹16,554"
$ws.Range("B8").Value = "Striado Rectangular Solid Wood Coffee Table In Teak Finish"
$ws.Range("C8").Value = "₹10,947"
$ws.Range("B9").Value = "Marcel Rectangular Metal Coffee Table In White Gloss Finish"
$ws.Range("C9").Value = "₹11,967"
$ws.Range("B10").Value = "Renesme Rectangular Solid Wood Coffee Table In Mahogany Finish"
$ws.Range("C10").Value = "₹15,317"
$ws.Range("B11").Value = "Dyson Abstract Metal Coffee Table In Teak Finish"
$ws.Range("C11").Value = "₹7,679"
$ws.Range("B12").Value = "Ivara Rectangular Solid Wood Coffee Table In Natural Finish"
$ws.Range("C12").Value = "₹16,049"
$ws.Range("B13").Value = "Botwin Rectangular Solid Wood Coffee Table In Mahogany Finish"
$ws.Range("C13").Value = "₹9,647"
$ws.Range("B14").Value = "Zephyr Rectangular Solid Wood Coffee Table In Teak Finish"
$ws.Range("C14").Value = "₹14,104"
$ws.Range("B15").Value = "Fring Engineered Wood Side Table In Matte Finish"
$ws.Range("C15").Value = "₹2,399"
$ws.Range("B16").Value = "Claire Rectangular Solid Wood Coffee Table In Mahogany Finish"
$ws.Range("C16").Value = "₹12,725"
$ws.Range("B17").Value = "Botwin Rectangular Solid Wood Coffee Table In Teak Finish"
$ws.Range("C17").Value = "₹9,647"
$ws.Range("B18").Value = "Epsilon Rectangular Solid Wood Coffee Table In Mahogany Finish"
$ws.Range("C18").Value = "₹11,384"
$ws.Range("B19").Value = "Dyson Rectangular Metal Coffee Table In Walnut Finish"
$ws.Range("C19").Value = "₹10,529"
$ws.Range("B20").Value = "Gustowe Rectangular Engineered Wood Coffee Table In Matte Finish"
$ws.Range("C20").Value = "₹2,279"
$ws.Range("B21").Value = "Striado Rectangular Solid Wood Coffee Table In Mahogany Finish"
$ws.Range("C21").Value = "₹10,947"
$ws.Range("B22").Value = "Osiris Rectangular Stone Coffee Table In Finish"
$ws.Range("C22").Value = "₹15,677"
$ws.Range("B23").Value = "Altura Rectangular Solid Wood Coffee Table In Two Tone Finish"
$ws.Range("C23").Value = "₹8,374"
$ws.Range("B24").Value = "Sylvie Rectangular Solid Wood Coffee Table In Natural Finish"
$ws.Range("C24").Value = "₹11,839"
$ws.Range("B25").Value = "Liam Rectangular Engineered Wood Coffee Table In Dark Wenge Finish"
$ws.Range("C25").Value = "₹3,817"
$ws.Range("B26").Value = "Florence Oval Solid Wood Coffee Table In Teak Finish"
$ws.Range("C26").Value = "₹10,223"
$ws.Range("B27").Value = "Reeves Rectangular Engineered Wood Coffee Table In Rustic Walnut Finish"
$ws.Range("C27").Value = "₹5,543"
$ws.Range("B28").Value = "Nitara Oval Solid Wood Coffee Table In Teak Finish"
$ws.Range("C28").Value = "₹12,095"
$ws.Range("B29").Value = "Renesme Rectangular Solid Wood Coffee Table In Teak Finish"
$ws.Range("C29").Value = "₹15,317"
$ws.Range("B30").Value = "Odette Square Solid Wood Coffee Table In Honey Oak Finish"
$ws.Range("C30").Value = "₹5,919"
$ws.Range("B31").Value = "Epsilon Rectangular Solid Wood Coffee Table In Teak Finish"
$ws.Range("C31").Value = "₹11,384"
